$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 22; existing rows 22:82 shift down to 23:83.
$ws.Rows("22").Insert()

# Populate the newly-inserted row 22 with the new data record
# (a weekly price entry for Espárragos at Feria Lagunitas de Puerto Montt).
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C22").Value = "Los Lagos"
$ws.Range("D22").Value = 45238
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 300000000
$ws.Range("G22").Value = "Espárragos"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 2000
$ws.Range("N22").Value = '$/kilo'
$ws.Range("O22").Value = "Provincia de Linares"
$ws.Range("P22").Value = 2000
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"
